$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "61.239.88"
$ws.Range("E2").Value = "  -0.02%  "

# Row 3
$ws.Range("D3").Value = "3.403.50"
$ws.Range("E3").Value = "  +0.97%  "

# Row 4
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.08%  "

# Row 5
$ws.Range("D5").Value = "'572.84"
$ws.Range("E5").Value = "  -0.12%  "

# Row 6
$ws.Range("D6").Value = "'138.91"
$ws.Range("E6").Value = "  +1.84%  "

# Row 7
$ws.Range("E7").Value = "  +0.10%  "

# Row 8
$ws.Range("D8").Value = "3.402.16"
$ws.Range("E8").Value = "  +0.97%  "

# Row 9
$ws.Range("D9").Value = "'0.469"
$ws.Range("E9").Value = "  -0.60%  "

# Row 10
$ws.Range("E10").Value = "  +2.86%  "

# Row 11
$ws.Range("E11").Value = "  -1.48%  "

# Row 12
$ws.Range("E12").Value = "  -1.91%  "

# Row 13
$ws.Range("D13").Value = "3.989.89"
$ws.Range("E13").Value = "  +1.05%  "

# Row 14
$ws.Range("E14").Value = "  -1.05%  "

# Row 15
$ws.Range("D15").Value = "'26.68"
$ws.Range("E15").Value = "  +2.53%  "

# Row 16
$ws.Range("E16").Value = "  -1.84%  "

# Row 17
$ws.Range("D17").Value = "3.399.79"
$ws.Range("E17").Value = "  +1.04%  "

# Row 18
$ws.Range("D18").Value = "61.328.21"
$ws.Range("E18").Value = "  -0.08%  "

# Row 19
$ws.Range("D19").Value = "'5.94"
$ws.Range("E19").Value = "  +1.59%  "

# Row 20
$ws.Range("E20").Value = "  -1.28%  "

# Row 21
$ws.Range("E21").Value = "  +0.48%  "

# Row 22
$ws.Range("D22").Value = "'378.17"
$ws.Range("E22").Value = "  +0.78%  "

# Row 23
$ws.Range("D23").Value = "3.528.27"
$ws.Range("E23").Value = "  +0.56%  "

# Row 24
$ws.Range("E24").Value = "  -0.28%  "

# Row 25
$ws.Range("E25").Value = "  +0.10%  "

# Row 26
$ws.Range("D26").Value = "'71.29"
$ws.Range("E26").Value = "  -0.28%  "

# Row 27
$ws.Range("E27").Value = "  -2.08%  "

# Row 28
$ws.Range("D28").Value = "'1.62"
$ws.Range("E28").Value = "  -4.30%  "

# Row 29
$ws.Range("D29").Value = "'0.174"
$ws.Range("E29").Value = "  +8.65%  "

# Row 30
$ws.Range("B30").Value = "Binance-PegBSC-USD"
$ws.Range("C30").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D30").Value = "'1.00"
$ws.Range("E30").Value = "  -0.29%  "

# Row 31
$ws.Range("B31").Value = "RenderToken"
$ws.Range("C31").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D31").Value = "'7.42"
$ws.Range("E31").Value = "  -0.34%  "

# Row 32
$ws.Range("D32").Value = "'8.10"
$ws.Range("E32").Value = "  -1.68%  "

# Row 33
$ws.Range("E33").Value = "  -0.93%  "

# Row 34
$ws.Range("E34").Value = "  -0.03%  "

# Row 35
$ws.Range("D35").Value = "'23.47"
$ws.Range("E35").Value = "  -0.16%  "

# Row 36
$ws.Range("E36").Value = "  +1.78%  "

# Row 37
$ws.Range("D37").Value = "'5.12"
$ws.Range("E37").Value = "  -2.79%  "

# Row 38
$ws.Range("D38").Value = "'6.86"
$ws.Range("E38").Value = "  +0.64%  "

# Row 39
$ws.Range("D39").Value = "'166.27"
$ws.Range("E39").Value = "  +0.29%  "

# Row 40
$ws.Range("D40").Value = "'0.0770"
$ws.Range("E40").Value = "  -0.50%  "

# Row 41
$ws.Range("D41").Value = "'26.28"
$ws.Range("E41").Value = "  +6.92%  "

# Row 42
$ws.Range("D42").Value = "'1.76"
$ws.Range("E42").Value = "  +2.70%  "

# Row 43
$ws.Range("D43").Value = "'0.999"
$ws.Range("E43").Value = "  -0.15%  "

# Row 44
$ws.Range("D44").Value = "'0.779"
$ws.Range("E44").Value = "  +0.78%  "

# Row 45
$ws.Range("D45").Value = "'41.96"
$ws.Range("E45").Value = "  +1.34%  "

# Row 46
$ws.Range("D46").Value = "'4.39"
$ws.Range("E46").Value = "  -0.36%  "

# Row 47
$ws.Range("D47").Value = "'1.17"
$ws.Range("E47").Value = "  -1.79%  "

# Row 48
$ws.Range("D48").Value = "2.530.19"
$ws.Range("E48").Value = "  +7.81%  "

# Row 49
$ws.Range("D49").Value = "'23.59"
$ws.Range("E49").Value = "  +4.46%  "

# Row 50
$ws.Range("E50").Value = "  -0.53%  "

# Row 51
$ws.Range("B51").Value = "VeChain"
$ws.Range("C51").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D51").Value = "'0.0262"
$ws.Range("E51").Value = "  +0.14%  "
